{"js": "// UC Visualizzazione Auto \u2014 \"correzione UC visualizzazione auto\"\n//\n// 1. Entry Condition cell: replace the whole sentence with a new one.\n// 2. \"2.a1\" label cell: collapse the two runs (\"2\" + \".a1\") into one run\n//    (bold formatting retained automatically since we rewrite in place).\n// 3. Error-message cell: collapse the three runs into a single run with\n//    the same (unchanged) text.\n// 4. \"2.a2\" label cell: collapse the two runs (\"2\" + \".a2\") into one run.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Helper: replace the full text of the first paragraph of a table cell,\n// using the paragraph's own Range so existing run/paragraph formatting\n// (bold, etc.) on that paragraph is preserved rather than reset.\nasync function replaceCellParagraphText(rowIndex, cellIndex, newText) {\n  const cell = table.getCell(rowIndex, cellIndex);\n  const paragraphs = cell.body.paragraphs;\n  paragraphs.load(\"items\");\n  await context.sync();\n\n  const paragraph = paragraphs.items[0];\n  const range = paragraph.getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) Entry Condition (row 7 / index 6, second cell of the merged row).\nawait replaceCellParagraphText(\n  6,\n  1,\n  \"L\\u2019utente \\u00e8 presente nell\\u2019home page.\"\n);\n\n// 2) \"2.a1\" label (row 20 / index 19, first cell).\nawait replaceCellParagraphText(19, 0, \"2.a1\");\n\n// 3) Error message (row 20 / index 19, third cell).\nawait replaceCellParagraphText(\n  19,\n  2,\n  \"Visualizza un messaggio di errore all\\u2019utente. Il messaggio segnala che il sistema non \\u00e8 riuscito ad effettuare il recupero dei dati.\"\n);\n\n// 4) \"2.a2\" label (row 21 / index 20, first cell).\nawait replaceCellParagraphText(20, 0, \"2.a2\");\n", "ps1": "# UC Visualizzazione Auto \u2014 \"correzione UC visualizzazione auto\"\n#\n# 1. Entry Condition cell: replace the whole sentence with a new one.\n# 2. \"2.a1\" label cell: collapse the two runs (\"2\" + \".a1\") into one run\n#    (bold formatting retained because Find/Replace rewrites in place).\n# 3. Error-message cell: collapse the three runs into a single run with\n#    the same (unchanged) text.\n# 4. \"2.a2\" label cell: collapse the two runs (\"2\" + \".a2\") into one run.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfunction Replace-CellParagraphText($table, $row, $col, $oldText, $newText) {\n    $cell = $table.Cell($row, $col)\n    $para = $cell.Range.Paragraphs.Item(1)\n    $rng = $para.Range\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.MatchWildcards = $false\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n}\n\n# 1) Entry Condition (table row 7, grid column 5 of the merged value cell).\nReplace-CellParagraphText $t 7 5 `\n    \"L\u2019utente deve essere loggato e clicca sull\u2019anteprima di un\u2019auto.\" `\n    \"L\u2019utente \u00e8 presente nell\u2019home page.\"\n\n# 2) \"2.a1\" label (table row 20, grid column 1).\nReplace-CellParagraphText $t 20 1 \"2.a1\" \"2.a1\"\n\n# 3) Error message (table row 20, grid column 5).\nReplace-CellParagraphText $t 20 5 `\n    \"Visualizza un messaggio di errore all\u2019utente. Il messaggio segnala che il sistema non \u00e8 riuscito ad effettuare il recupero dei dati.\" `\n    \"Visualizza un messaggio di errore all\u2019utente. Il messaggio segnala che il sistema non \u00e8 riuscito ad effettuare il recupero dei dati.\"\n\n# 4) \"2.a2\" label (table row 21, grid column 1).\nReplace-CellParagraphText $t 21 1 \"2.a2\" \"2.a2\"\n"}
